$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 137
$ws.Range("H137").Value2 = 1105.1666   # was 966.6774
$ws.Range("I137").Value2 = 1040.5454   # was 876.8182
$ws.Range("J137").Value2 = 1206.7142   # was 1186.3334
$ws.Range("K137").Value2 = 3121.6362   # was 2630.4546
$ws.Range("L137").Value2 = 3620.1426   # was 3559.0002
$ws.Range("M137").Value2 = -571.6361999999999   # was -80.45460000000003
$ws.Range("N137").Value2 = -8720.142599999999   # was -8659.0002
# Row 138
$ws.Range("H138").Value2 = 2124.4856   # was 2205.5151
$ws.Range("I138").Value2 = 1537.1111   # was 1630.8125
$ws.Range("K138").Value2 = 4611.3333   # was 4892.4375
$ws.Range("M138").Value2 = 528.6666999999998   # was 247.5625

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value2 = 6727.5874   # was 6536.7935
$ws.Range("I32").Value2 = 2720.328   # was 2523.2788
$ws.Range("K32").Value2 = 2720.328   # was 2523.2788
$ws.Range("M32").Value2 = -2433.328   # was -2236.2788
# Row 61
$ws.Range("H61").Value2 = 2055.3076   # was 2117.3333
$ws.Range("I61").Value2 = 1996.1818   # was 2064.7
$ws.Range("K61").Value2 = 1996.1818   # was 2064.7
$ws.Range("M61").Value2 = -1784.1818   # was -1852.7
# Row 63
$ws.Range("H63").Value2 = 2949   # was 1427.8
$ws.Range("I63").Value2 = 2949   # was 1427.8
$ws.Range("K63").Value2 = 2949   # was 1427.8
$ws.Range("M63").Value2 = -2263   # was -741.8
# Row 66
$ws.Range("H66").Value2 = 2949   # was 1427.8
$ws.Range("I66").Value2 = 2949   # was 1427.8
$ws.Range("K66").Value2 = 14745   # was 7139
$ws.Range("M66").Value2 = -11313   # was -3707
# Row 74
$ws.Range("H74").Value2 = 1320.1578   # was 1382.8334
$ws.Range("I74").Value2 = 1210.7646   # was 1274.4375
$ws.Range("K74").Value2 = 1210.7646   # was 1274.4375
$ws.Range("M74").Value2 = -336.7646   # was -400.4375
# Row 77
$ws.Range("H77").Value2 = 1320.1578   # was 1382.8334
$ws.Range("I77").Value2 = 1210.7646   # was 1274.4375
$ws.Range("K77").Value2 = 6053.823   # was 6372.1875
$ws.Range("M77").Value2 = -1685.823   # was -2004.1875
# Row 88
$ws.Range("H88").Value2 = 14286535   # was 16667491
$ws.Range("I88").Value2 = 800.6667   # was 768.3333
$ws.Range("J88").Value2 = 25000836   # was 22223066
$ws.Range("K88").Value2 = 800.6667   # was 768.3333
$ws.Range("L88").Value2 = 25000836   # was 22223066
$ws.Range("M88").Value2 = -394.6667   # was -362.3333
$ws.Range("N88").Value2 = -25001648   # was -22223878
# Row 91
$ws.Range("H91").Value2 = 14286535   # was 16667491
$ws.Range("I91").Value2 = 800.6667   # was 768.3333
$ws.Range("J91").Value2 = 25000836   # was 22223066
$ws.Range("K91").Value2 = 800.6667   # was 768.3333
$ws.Range("L91").Value2 = 25000836   # was 22223066
$ws.Range("M91").Value2 = 603.3333   # was 635.6667
$ws.Range("N91").Value2 = -25003644   # was -22225874
# Row 132
$ws.Range("H132").Value2 = 1597.0869   # was 1563.875
$ws.Range("I132").Value2 = 1535.6364   # was 1503.6522
$ws.Range("K132").Value2 = 4606.9092   # was 4510.9566
$ws.Range("M132").Value2 = -2076.9092   # was -1980.9566
# Row 135
$ws.Range("H135").Value2 = 61901.168   # was 61898.832
$ws.Range("J135").Value2 = 61901.168   # was 61898.832
$ws.Range("L135").Value2 = 61901.168   # was 61898.832
$ws.Range("N135").Value2 = -72041.16800000001   # was -72038.83199999999
# Row 136
$ws.Range("H136").Value2 = 2055.3076   # was 2117.3333
$ws.Range("I136").Value2 = 1996.1818   # was 2064.7
$ws.Range("K136").Value2 = 5988.5454   # was 6194.099999999999
$ws.Range("M136").Value2 = -3438.5454   # was -3644.099999999999

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value2 = 2465.0588   # was 2488.4707
$ws.Range("I86").Value2 = 2336.818   # was 2373
$ws.Range("K86").Value2 = 2336.818   # was 2373
$ws.Range("M86").Value2 = -1213.818   # was -1250
# Row 89
$ws.Range("H89").Value2 = 2465.0588   # was 2488.4707
$ws.Range("I89").Value2 = 2336.818   # was 2373
$ws.Range("K89").Value2 = 11684.09   # was 11865
$ws.Range("M89").Value2 = -6068.09   # was -6249
# Row 137
$ws.Range("H137").Value2 = 50000   # was 0
$ws.Range("J137").Value2 = 50000   # was 0
$ws.Range("L137").Value2 = 50000   # was 0
$ws.Range("N137").Value2 = -60200

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value2 = 13209.637   # was 13958.6455
$ws.Range("I31").Value2 = 3311.5652   # was 3407.5454
$ws.Range("J31").Value2 = 35975.2   # was 39750.223
$ws.Range("K31").Value2 = 3311.5652   # was 3407.5454
$ws.Range("L31").Value2 = 35975.2   # was 39750.223
$ws.Range("M31").Value2 = -3016.5652   # was -3112.5454
$ws.Range("N31").Value2 = -36565.2   # was -40340.223
# Row 34
$ws.Range("H34").Value2 = 13209.637   # was 13958.6455
$ws.Range("I34").Value2 = 3311.5652   # was 3407.5454
$ws.Range("J34").Value2 = 35975.2   # was 39750.223
$ws.Range("K34").Value2 = 3311.5652   # was 3407.5454
$ws.Range("L34").Value2 = 35975.2   # was 39750.223
$ws.Range("M34").Value2 = -3109.5652   # was -3205.5454
$ws.Range("N34").Value2 = -36379.2   # was -40154.223
# Row 97
$ws.Range("H97").Value2 = 35764.855   # was 34669.25
$ws.Range("J97").Value2 = 35764.855   # was 34669.25
$ws.Range("L97").Value2 = 35764.855   # was 34669.25
$ws.Range("N97").Value2 = -37746.855   # was -36651.25
# Row 99
$ws.Range("H99").Value2 = 9468.708000000001   # was 9706.478999999999
$ws.Range("J99").Value2 = 4422.625   # was 4483
$ws.Range("L99").Value2 = 4422.625   # was 4483
$ws.Range("N99").Value2 = -7418.625   # was -7479
# Row 126
$ws.Range("H126").Value2 = 9468.708000000001   # was 9706.478999999999
$ws.Range("J126").Value2 = 4422.625   # was 4483
$ws.Range("L126").Value2 = 13267.875   # was 13449
$ws.Range("N126").Value2 = -18207.875   # was -18389
# Row 132
$ws.Range("H132").Value2 = 3200.125   # was 2925.75
$ws.Range("I132").Value2 = 3086.5862   # was 2801.0303
$ws.Range("K132").Value2 = 9259.758600000001   # was 8403.090899999999
$ws.Range("M132").Value2 = -6729.758600000001   # was -5873.090899999999
# Row 134
$ws.Range("H134").Value2 = 2428.913   # was 2304.923
$ws.Range("I134").Value2 = 1486.3889   # was 1467.5238
$ws.Range("K134").Value2 = 4459.1667   # was 4402.5714
$ws.Range("M134").Value2 = -1924.1667   # was -1867.5714

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 16
$ws.Range("H16").Value2 = 0   # was 346.66666
$ws.Range("I16").Value2 = 0   # was 346.66666
$ws.Range("K16").Value2 = 0   # was 1039.99998
$ws.Range("M16").ClearContents()
# Row 29
$ws.Range("H29").Value2 = 89.8   # was 100
$ws.Range("I29").Value2 = 124.5   # was 200
$ws.Range("K29").Value2 = 373.5   # was 600
$ws.Range("M29").Value2 = -96.5   # was -323
# Row 34
$ws.Range("H34").Value2 = 651.6429000000001   # was 718.9375
$ws.Range("J34").Value2 = 1458   # was 1381.4286
$ws.Range("L34").Value2 = 4374   # was 4144.2858
$ws.Range("N34").Value2 = -4542   # was -4312.2858
# Row 39
$ws.Range("H39").Value2 = 4291.3335   # was 3833.8
$ws.Range("J39").Value2 = 4291.3335   # was 3833.8
$ws.Range("L39").Value2 = 12874.0005   # was 11501.4
$ws.Range("N39").Value2 = -13462.0005   # was -12089.4
# Row 55
$ws.Range("H55").Value2 = 8930678   # was 8335578.5
$ws.Range("I55").Value2 = 699.6   # was 749.75
$ws.Range("J55").Value2 = 13891776   # was 11366425
$ws.Range("K55").Value2 = 2098.8   # was 2249.25
$ws.Range("L55").Value2 = 41675328   # was 34099275
$ws.Range("M55").Value2 = -1921.8   # was -2072.25
$ws.Range("N55").Value2 = -41675682   # was -34099629
# Row 101
$ws.Range("H101").Value2 = 14999   # was 9998
$ws.Range("J101").Value2 = 20000   # was 0
$ws.Range("L101").Value2 = 60000   # was 0
$ws.Range("N101").Value2 = -64868
# Row 114
$ws.Range("H114").Value2 = 40002340   # was 66667230
$ws.Range("J114").Value2 = 5000   # was 0
$ws.Range("L114").Value2 = 15000   # was 0
$ws.Range("N114").Value2 = -21508
# Row 129
$ws.Range("H129").Value2 = 3709.4   # was 3228.6667
$ws.Range("J129").Value2 = 6299   # was 4735
$ws.Range("L129").Value2 = 18897   # was 14205
$ws.Range("N129").Value2 = -28897   # was -24205
# Row 140
$ws.Range("H140").Value2 = 1875.5714   # was 1627.8182
$ws.Range("I140").Value2 = 1875.5714   # was 1577.3
$ws.Range("J140").Value2 = 0   # was 2133
$ws.Range("K140").Value2 = 5626.7142   # was 4731.9
$ws.Range("L140").Value2 = 0   # was 6399
$ws.Range("M140").Value2 = -446.7142000000003   # was 448.1000000000004
$ws.Range("N140").ClearContents()

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value2 = 41529   # was 42893.324
$ws.Range("I97").Value2 = 27749.75   # was 28989.521
$ws.Range("K97").Value2 = 27749.75   # was 28989.521
$ws.Range("M97").Value2 = -27253.75   # was -28493.521
# Row 124
$ws.Range("H124").Value2 = 0   # was 151080
$ws.Range("J124").Value2 = 0   # was 151080
$ws.Range("L124").Value2 = 0   # was 151080
$ws.Range("N124").ClearContents()
# Row 133
$ws.Range("H133").Value2 = 98765   # was 0
$ws.Range("J133").Value2 = 98765   # was 0
$ws.Range("L133").Value2 = 98765   # was 0
$ws.Range("N133").Value2 = -108885

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value2 = 15908.259   # was 17342.36
$ws.Range("I7").Value2 = 17996.592   # was 19887.8
$ws.Range("J7").Value2 = 6719.6   # was 7160.6
$ws.Range("K7").Value2 = 17996.592   # was 19887.8
$ws.Range("L7").Value2 = 6719.6   # was 7160.6
$ws.Range("M7").Value2 = -17884.592   # was -19775.8
$ws.Range("N7").Value2 = -6943.6   # was -7384.6
# Row 99
$ws.Range("H99").Value2 = 25000   # was 30285
$ws.Range("I99").Value2 = 25000   # was 0
$ws.Range("J99").Value2 = 0   # was 30285
$ws.Range("K99").Value2 = 25000   # was 0
$ws.Range("L99").Value2 = 0   # was 30285
$ws.Range("N99").ClearContents()
$ws.Range("M99").Value2 = -22005
# Row 126
$ws.Range("H126").Value2 = 15908.259   # was 17342.36
$ws.Range("I126").Value2 = 17996.592   # was 19887.8
$ws.Range("J126").Value2 = 6719.6   # was 7160.6
$ws.Range("K126").Value2 = 53989.776   # was 59663.39999999999
$ws.Range("L126").Value2 = 20158.8   # was 21481.8
$ws.Range("M126").Value2 = -51519.776   # was -57193.39999999999
$ws.Range("N126").Value2 = -25098.8   # was -26421.8
# Row 132
$ws.Range("H132").Value2 = 3845.2222   # was 3818.5715
$ws.Range("I132").Value2 = 3689.8635   # was 3664.1738
$ws.Range("K132").Value2 = 11069.5905   # was 10992.5214
$ws.Range("M132").Value2 = -8539.5905   # was -8462.5214
# Row 133
$ws.Range("H133").Value2 = 109998.75   # was 113332
$ws.Range("J133").Value2 = 109998.75   # was 113332
$ws.Range("L133").Value2 = 109998.75   # was 113332
$ws.Range("N133").Value2 = -115058.75   # was -118392
# Row 136
$ws.Range("H136").Value2 = 3677.394   # was 3294.0264
$ws.Range("I136").Value2 = 3157.652   # was 2730.1785
$ws.Range("K136").Value2 = 9472.956   # was 8190.5355
$ws.Range("M136").Value2 = -6922.956   # was -5640.5355

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value2 = 3032.524   # was 2857.9565
$ws.Range("I100").Value2 = 3345.923   # was 3146.2144
$ws.Range("J100").Value2 = 2523.25   # was 2409.5557
$ws.Range("K100").Value2 = 6691.846   # was 6292.4288
$ws.Range("L100").Value2 = 5046.5   # was 4819.1114
$ws.Range("M100").Value2 = -6150.846   # was -5751.4288
$ws.Range("N100").Value2 = -6128.5   # was -5901.1114
# Row 127
$ws.Range("H127").Value2 = 131657.33   # was 129992.664
$ws.Range("J127").Value2 = 131657.33   # was 129992.664
$ws.Range("L127").Value2 = 131657.33   # was 129992.664
$ws.Range("N127").Value2 = -141577.33   # was -139912.664
# Row 136
$ws.Range("H136").Value2 = 2967.68   # was 3129.2173
$ws.Range("I136").Value2 = 2611.1   # was 2777.889
$ws.Range("K136").Value2 = 7833.299999999999   # was 8333.667000000001
$ws.Range("M136").Value2 = -5283.299999999999   # was -5783.667000000001
